# Primer version funcional con GUI
# Reemplazo de nombres de columnas (Factura/Diploma -> rutas completas) y
# valores de archivos adjuntos (agregado de extensión .png), mas wrap text
# en la columna de Diplomas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$facturasPath = "C:\Users\Lucas\Desktop\Coding\Github\NEXXOS-emails\Facturas"
$diplomasPath = "C:\Users\Lucas\Desktop\Coding\Github\NEXXOS-emails\Diplomas"

# Update header names for columns D (Factura) and E (Diploma)
$ws.Range("D1").Value = $facturasPath
$ws.Range("E1").Value = $diplomasPath

# Update attachment filenames to include the .png extension
$ws.Range("D2").Value = "abbate-f.png"
$ws.Range("E2").Value = "abbate-d.png"
$ws.Range("D3").Value = "capriata-f.png"
$ws.Range("E3").Value = "capriata-d.png"

# Widen column E and slightly narrow column D
$ws.Columns.Item(4).ColumnWidth = 21.666666666666668
$ws.Columns.Item(5).ColumnWidth = 46.666666666666664

# Wrap text only for the cells that actually hold (or will hold) data in
# column E, mirroring the author's per-cell formatting rather than a blanket
# column-level format.
$ws.Range("E1").WrapText = $true
$ws.Range("E2").WrapText = $true
$ws.Range("E3").WrapText = $true
$ws.Range("E11").WrapText = $true

# Make the header row taller to accommodate the wrapped header text
$ws.Rows.Item(1).RowHeight = 30

# Move active selection to E19 (matches author's cursor position on save)
$ws.Range("E19").Select()
